# Apply Betfair Back/Lay odds updates (2025-10-09) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.47
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 2.7
$ws.Range("I2").Value = 150
$ws.Range("J2").Value = 3.45
$ws.Range("K2").Value = 950
$ws.Range("L2").Value = 1.25
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 1.11
$ws.Range("P2").Value = 1.89
$ws.Range("R2").Value = 1.18
$ws.Range("S2").Value = 2.04
$ws.Range("T2").Value = 1.04
$ws.Range("V2").Value = 1.23
$ws.Range("W2").Value = 1.4
$ws.Range("X2").Value = 60

# Row 3
$ws.Range("F3").Value = 1.26
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 2.92
$ws.Range("J3").Value = 3.05
$ws.Range("K3").Value = 980
$ws.Range("L3").Value = 1.25
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 1.73
$ws.Range("R3").Value = 1.28
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.03
$ws.Range("V3").Value = 1.07
$ws.Range("W3").Value = 2
$ws.Range("AN3").Value = 600

# Row 4
$ws.Range("F4").Value = 2.64
$ws.Range("G4").Value = 2.74
$ws.Range("H4").Value = 2.96
$ws.Range("I4").Value = 3.1
$ws.Range("K4").Value = 3.35
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 2.98
$ws.Range("O4").Value = 1.42
$ws.Range("P4").Value = 1.68
$ws.Range("S4").Value = 4.4
$ws.Range("V4").Value = 1.47
$ws.Range("AA4").Value = 130
$ws.Range("AE4").Value = 190
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 110
$ws.Range("AK4").Value = 85
$ws.Range("AM4").Value = 580

# Row 5
$ws.Range("F5").Value = 1.86
$ws.Range("G5").Value = 1.96
$ws.Range("H5").Value = 5.2
$ws.Range("I5").Value = 5.9
$ws.Range("K5").Value = 3.55
$ws.Range("N5").Value = 2.6
$ws.Range("O5").Value = 1.54
$ws.Range("P5").Value = 1.53
$ws.Range("Q5").Value = 2.6
$ws.Range("R5").Value = 1.18
$ws.Range("S5").Value = 5.6
$ws.Range("T5").Value = 2.26
$ws.Range("U5").Value = 1.68
$ws.Range("X5").Value = 9
$ws.Range("Y5").Value = 36
$ws.Range("Z5").Value = 980
$ws.Range("AD5").Value = 80
$ws.Range("AF5").Value = 19.5
$ws.Range("AG5").Value = 40
$ws.Range("AI5").Value = 540
$ws.Range("AJ5").Value = 900
$ws.Range("AK5").Value = 110
$ws.Range("AL5").Value = 450
$ws.Range("AN5").Value = 980
$ws.Range("AO5").Value = 230

# Row 6
$ws.Range("G6").Value = 1.95
$ws.Range("I6").Value = 5.3
$ws.Range("W6").Value = 2.04

# Row 7
$ws.Range("F7").Value = 3.1
$ws.Range("G7").Value = 3.45
$ws.Range("K7").Value = 3.8
$ws.Range("S7").Value = 3.75
$ws.Range("T7").Value = 1.78
$ws.Range("U7").Value = 2.04
$ws.Range("Y7").Value = 10.5
$ws.Range("AA7").Value = 980
$ws.Range("AK7").Value = 980
$ws.Range("AL7").Value = 980
$ws.Range("AN7").Value = 980

# Row 8
$ws.Range("F8").Value = 3.6
$ws.Range("G8").Value = 3.9
$ws.Range("H8").Value = 2.24
$ws.Range("I8").Value = 2.42
$ws.Range("J8").Value = 3.15
$ws.Range("K8").Value = 3.45
$ws.Range("N8").Value = 2.8
$ws.Range("P8").Value = 1.61
$ws.Range("V8").Value = 1.7
$ws.Range("W8").Value = 1.34
$ws.Range("AA8").Value = 1000
$ws.Range("AC8").Value = 7.6
$ws.Range("AG8").Value = 1000

# Row 9
$ws.Range("H9").Value = 2.22
$ws.Range("J9").Value = 1.09
$ws.Range("K9").Value = 5.2
$ws.Range("L9").Value = 1.45
$ws.Range("N9").Value = 1.3
$ws.Range("P9").Value = 1.3
$ws.Range("V9").Value = 1.5
